$d = $word.ActiveDocument

# The sentence currently reads "Abdulrahman is quicker than Ameer."
# and needs to become "Abdulrahman is faster than Ameer." - but,
# per the target markup, split across three runs:
#   "Abdulrahman is " | "faster" | " than Ameer."
$oldWord = "quicker"
$newWord = "faster"

# Find the word to swap out.
$findRange = $d.Content
$found = $findRange.Find.Execute($oldWord, $true, $false, $false, $false,
                                  $false, $true, 1, $false, "", 0)

$start = $findRange.Start

# Replace the text in place (this keeps everything as a single run for
# the moment - Word / this host coalesces same-formatted runs back
# together on save).
$findRange.Text = $newWord
$end = $start + $newWord.Length

# Re-select exactly the newly-inserted replacement word and wrap it in a
# throwaway bookmark. Adding (and then immediately deleting) a bookmark
# forces a run split at the bookmark's boundaries, which is what gives
# us the three separate <w:r> runs the edit calls for - without that,
# the host merges adjacent runs that share identical run formatting
# right back into one run.
$splitRange = $d.Range($start, $end)
$markName = "__tmp_split__"
$d.Bookmarks.Add($markName, $splitRange)
$d.Bookmarks($markName).Delete()
